# Auto-generated edit script: updates LeveProfit sheet price/profit columns (H,I,J,K,L,M,N)
# per scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3928.625
$ws.Range("I76").Value = 3989.8572
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 3989.8572
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -3674.8572
$ws.Range("N76").Value = -4130

$ws.Range("H79").Value = 3928.625
$ws.Range("I79").Value = 3989.8572
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 3989.8572
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -2897.8572
$ws.Range("N79").Value = -5684

$ws.Range("H138").Value = 3138.5938
$ws.Range("J138").Value = 4584.8237
$ws.Range("L138").Value = 13754.4711
$ws.Range("N138").Value = -24034.4711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 868.875
$ws.Range("I5").Value = 942
$ws.Range("K5").Value = 942
$ws.Range("M5").Value = -830

$ws.Range("H32").Value = 20866.375
$ws.Range("I32").Value = 18602.885
$ws.Range("K32").Value = 18602.885
$ws.Range("M32").Value = -18315.885

$ws.Range("H40").Value = 34941.117
$ws.Range("J40").Value = 34756.25
$ws.Range("L40").Value = 34756.25
$ws.Range("N40").Value = -35108.25

$ws.Range("H63").Value = 3429.2068
$ws.Range("I63").Value = 2464.6191
$ws.Range("J63").Value = 5961.25
$ws.Range("K63").Value = 2464.6191
$ws.Range("L63").Value = 5961.25
$ws.Range("M63").Value = -1778.6191
$ws.Range("N63").Value = -7333.25

$ws.Range("H66").Value = 3429.2068
$ws.Range("I66").Value = 2464.6191
$ws.Range("J66").Value = 5961.25
$ws.Range("K66").Value = 12323.0955
$ws.Range("L66").Value = 29806.25
$ws.Range("M66").Value = -8891.0955
$ws.Range("N66").Value = -36670.25

$ws.Range("H88").Value = 1303.1875
$ws.Range("J88").Value = 1752.625
$ws.Range("L88").Value = 1752.625
$ws.Range("N88").Value = -2564.625

$ws.Range("H91").Value = 1303.1875
$ws.Range("J91").Value = 1752.625
$ws.Range("L91").Value = 1752.625
$ws.Range("N91").Value = -4560.625

$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -109799

$ws.Range("H132").Value = 2605.673
$ws.Range("I132").Value = 1996.2128
$ws.Range("K132").Value = 5988.6384
$ws.Range("M132").Value = -3458.6384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 868.875
$ws.Range("I4").Value = 942
$ws.Range("K4").Value = 942
$ws.Range("M4").Value = -827

$ws.Range("H107").Value = 6640.717
$ws.Range("I107").Value = 7170.0576
$ws.Range("K107").Value = 7170.0576
$ws.Range("M107").Value = -5250.0576

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 309.33334
$ws.Range("J7").Value = 494
$ws.Range("L7").Value = 494
$ws.Range("N7").Value = -720

$ws.Range("H58").Value = 17584.285
$ws.Range("I58").Value = 1397.5
$ws.Range("K58").Value = 1397.5
$ws.Range("M58").Value = -1194.5

$ws.Range("H136").Value = 17584.285
$ws.Range("I136").Value = 1397.5
$ws.Range("K136").Value = 4192.5
$ws.Range("M136").Value = -1642.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1064.625
$ws.Range("J98").Value = 1002.75
$ws.Range("L98").Value = 3008.25
$ws.Range("N98").Value = -6004.25

$ws.Range("H117").Value = 1384.4546
$ws.Range("I117").Value = 1042.6666
$ws.Range("K117").Value = 3127.9998
$ws.Range("M117").Value = 314.0001999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5661.75
$ws.Range("I70").Value = 6201.75
$ws.Range("J70").Value = 5121.75
$ws.Range("K70").Value = 6201.75
$ws.Range("L70").Value = 5121.75
$ws.Range("M70").Value = -5931.75
$ws.Range("N70").Value = -5661.75

$ws.Range("H73").Value = 5661.75
$ws.Range("I73").Value = 6201.75
$ws.Range("J73").Value = 5121.75
$ws.Range("K73").Value = 6201.75
$ws.Range("L73").Value = 5121.75
$ws.Range("M73").Value = -5265.75
$ws.Range("N73").Value = -6993.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 64999.75
$ws.Range("J38").Value = 64999.75
$ws.Range("L38").Value = 64999.75
$ws.Range("N38").Value = -65819.75

$ws.Range("H55").Value = 1056.0834
$ws.Range("I55").Value = 630.4
$ws.Range("K55").Value = 630.4
$ws.Range("M55").Value = -457.4

$ws.Range("H68").Value = 3731.7144
$ws.Range("I68").Value = 1522.9375
$ws.Range("J68").Value = 10799.8
$ws.Range("K68").Value = 1522.9375
$ws.Range("L68").Value = 10799.8
$ws.Range("M68").Value = -773.9375
$ws.Range("N68").Value = -12297.8

$ws.Range("H71").Value = 3731.7144
$ws.Range("I71").Value = 1522.9375
$ws.Range("J71").Value = 10799.8
$ws.Range("K71").Value = 7614.6875
$ws.Range("L71").Value = 53999
$ws.Range("M71").Value = -3870.6875
$ws.Range("N71").Value = -61487

$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62246

$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -191232

$ws.Range("H132").Value = 3508.362
$ws.Range("I132").Value = 3213.4348
$ws.Range("J132").Value = 4638.9165
$ws.Range("K132").Value = 9640.3044
$ws.Range("L132").Value = 13916.7495
$ws.Range("M132").Value = -7110.304400000001
$ws.Range("N132").Value = -18976.7495

$ws.Range("H133").Value = 99511
$ws.Range("J133").Value = 99511
$ws.Range("L133").Value = 99511
$ws.Range("N133").Value = -104571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 12002459
$ws.Range("J11").Value = 13336266
$ws.Range("L11").Value = 13336266
$ws.Range("N11").Value = -13336550

$ws.Range("H69").Value = 107527.25
$ws.Range("J69").Value = 107527.25
$ws.Range("L69").Value = 107527.25
$ws.Range("N69").Value = -109025.25

$ws.Range("H72").Value = 107527.25
$ws.Range("J72").Value = 107527.25
$ws.Range("L72").Value = 322581.75
$ws.Range("N72").Value = -330069.75

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
